$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures ---
# VALOR MORA total
$ws.Range("E11").Value = 88000
# Cant. Trabajadores (worker count)
$ws.Range("C13").Value = 1

# --- Update the remaining worker record (row 16) with the new data ---
$ws.Range("C16").Value = "1047367908"
$ws.Range("D16").Value = "EDWIN ANTONIO MENCO VANEGAS"
$ws.Range("E16").Value = "2508"
$ws.Range("F16").Value = 88000
$ws.Range("G16").Value = 2200000

# --- Remove the other three worker records (rows 17-19) ---
# This shifts everything below (including the signature block) up by 3 rows.
$ws.Range("A17:A19").EntireRow.Delete()
